# Auto-generated edit script to update '想去人数' (F column) values
# per commit: Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 95
$ws.Cells.Item(5, 6).Value = 9461
$ws.Cells.Item(6, 6).Value = 7286
$ws.Cells.Item(7, 6).Value = 7875
$ws.Cells.Item(10, 6).Value = 42
$ws.Cells.Item(11, 6).Value = 6782
$ws.Cells.Item(13, 6).Value = 490
$ws.Cells.Item(14, 6).Value = 462
$ws.Cells.Item(16, 6).Value = 667
$ws.Cells.Item(19, 6).Value = 225
$ws.Cells.Item(22, 6).Value = 116
$ws.Cells.Item(23, 6).Value = 11014
$ws.Cells.Item(24, 6).Value = 96
$ws.Cells.Item(25, 6).Value = 65
$ws.Cells.Item(26, 6).Value = 2089
$ws.Cells.Item(27, 6).Value = 2747
$ws.Cells.Item(29, 6).Value = 2464
$ws.Cells.Item(34, 6).Value = 2246
$ws.Cells.Item(36, 6).Value = 1527
$ws.Cells.Item(38, 6).Value = 50
$ws.Cells.Item(39, 6).Value = 5573
$ws.Cells.Item(42, 6).Value = 796
$ws.Cells.Item(43, 6).Value = 146
$ws.Cells.Item(44, 6).Value = 181
$ws.Cells.Item(45, 6).Value = 1094
$ws.Cells.Item(46, 6).Value = 1037
$ws.Cells.Item(47, 6).Value = 1458
$ws.Cells.Item(48, 6).Value = 83
$ws.Cells.Item(49, 6).Value = 1117
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 4
$ws.Cells.Item(16, 6).Value = 105
$ws.Cells.Item(21, 6).Value = 10
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 137
$ws.Cells.Item(3, 6).Value = 238
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 95
$ws.Cells.Item(4, 6).Value = 9461
$ws.Cells.Item(5, 6).Value = 9461
$ws.Cells.Item(6, 6).Value = 7286
$ws.Cells.Item(7, 6).Value = 137
$ws.Cells.Item(8, 6).Value = 238
$ws.Cells.Item(11, 6).Value = 7875
$ws.Cells.Item(14, 6).Value = 42
$ws.Cells.Item(15, 6).Value = 6782
$ws.Cells.Item(16, 6).Value = 6782
$ws.Cells.Item(18, 6).Value = 490
$ws.Cells.Item(19, 6).Value = 462
$ws.Cells.Item(20, 6).Value = 667
$ws.Cells.Item(23, 6).Value = 225
$ws.Cells.Item(28, 6).Value = 11014
$ws.Cells.Item(29, 6).Value = 96
$ws.Cells.Item(30, 6).Value = 65
$ws.Cells.Item(31, 6).Value = 2089
$ws.Cells.Item(32, 6).Value = 2747
$ws.Cells.Item(33, 6).Value = 2464
$ws.Cells.Item(37, 6).Value = 2246
$ws.Cells.Item(39, 6).Value = 1527
$ws.Cells.Item(40, 6).Value = 50
$ws.Cells.Item(41, 6).Value = 5573
$ws.Cells.Item(44, 6).Value = 796
$ws.Cells.Item(45, 6).Value = 146
$ws.Cells.Item(46, 6).Value = 181
$ws.Cells.Item(47, 6).Value = 1094
$ws.Cells.Item(48, 6).Value = 1037
$ws.Cells.Item(49, 6).Value = 1458
$ws.Cells.Item(50, 6).Value = 83
$ws.Cells.Item(51, 6).Value = 1117
